$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.170.93'
$ws.Range("E2").Value = '  +2.75%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.808.66'
$ws.Range("E3").Value = '  +0.75%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.49'
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3932'
$ws.Range("E7").Value = '  +3.53%  '
$ws.Range("E8").Value = '  +0.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.33'
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07544'
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.08'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.507'
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.811.48'
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.147'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06698'
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.83'
$ws.Range("E19").Value = '  -0.11%  '
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.75'
$ws.Range("E21").Value = '  +2.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.579'
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.156.74'
$ws.Range("E23").Value = '  +2.68%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.44'
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.412'
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.488'
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.530'
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.33'
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '153.56'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.017.74'
$ws.Range("E30").Value = '  +0.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '135.37'
$ws.Range("E31").Value = '  +1.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.164'
$ws.Range("E32").Value = '  +1.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.024'
$ws.Range("E33").Value = '  -1.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08851'
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.05'
$ws.Range("E35").Value = '  -1.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6937'
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06553'
$ws.Range("E37").Value = '  +2.56%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.451'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02418'
$ws.Range("E39").Value = '  +2.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.609'
$ws.Range("E40").Value = '  -2.76%  '
$ws.Range("E41").Value = '  +0.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.257'
$ws.Range("E42").Value = '  -1.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.466'
$ws.Range("E43").Value = '  -4.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.62'
$ws.Range("E44").Value = '  +1.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6418'
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.874'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.142'
$ws.Range("E48").Value = '  +0.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '131.39'
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07194'
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.92'
$ws.Range("E51").Value = '  +0.43%  '
